$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# A new stock item ("OFLAM 50MG 20 QUICKETAB.") needs to be inserted into the
# sorted product list as row 15 (between "NORHINOSE ..." on row 14 and
# "PAROFEN ..." which is currently on row 15). Inserting the row pushes every
# row from 15 downward down by one (old row 15 -> new row 16, ... old row 24
# -> new row 25).
# ---------------------------------------------------------------------------

$ws.Rows.Item(15).Insert()

# Re-create the per-column formatting for the new row by copying it from the
# row above (row 14, "NORHINOSE ..."), which carries the same style set used
# by every product row in the table.
for ($col = 1; $col -le 17; $col++) {
    $ws.Cells.Item(14, $col).Copy($ws.Cells.Item(15, $col))
}

# Match the row height used by the other product rows.
$ws.Rows.Item(15).RowHeight = 24.75

# Re-create the merged cells for the new product row (matches the pattern
# used by every other product row, e.g. row 14/16).
$ws.Range("A15:B15").Merge()
$ws.Range("C15:G15").Merge()
$ws.Range("H15:K15").Merge()
$ws.Range("L15:M15").Merge()
$ws.Range("N15:O15").Merge()

# Fill in the new row's values.
$ws.Range("A15").Value = 9
$ws.Range("C15").Value = "OFLAM 50MG 20 QUICKETAB."
$ws.Range("H15").Value = "0:1"
$ws.Range("L15").Value = "1"
$ws.Range("N15").Value = "48.00"
$ws.Range("P15").Value = "24.0000"
$ws.Range("Q15").Value = "0:1"

# ---------------------------------------------------------------------------
# The grand-total cell (previously P23) now lives at P24; bump it by the
# newly inserted item's sell price (1319.02 + 24.00 = 1343.02).
# ---------------------------------------------------------------------------
$ws.Range("P24").Value = 1343.02

# ---------------------------------------------------------------------------
# Refresh the "printed at" timestamp (now on row 25) to the new export time.
# ---------------------------------------------------------------------------
$ws.Range("A25").Value = "Friday, 15 August, 2025 7:39 PM"
